$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row 64 entries
$ws.Range("B64").Value = "Réalisation"
$ws.Range("C64").Value = 2.5
$ws.Range("D64").Value = "Création de la page d'affichage pour un article"
$ws.Range("E64").Value = "Accessible (pour le moment) que depuis le carousel de la home page"

# Match row 63's taller row height (content now wraps across two lines)
$ws.Rows.Item(64).RowHeight = 30

# Update the selected cell in the sheet view
$ws.Range("C65").Select()
